# Updates the cryptos price list with refreshed price/volume figures and
# restores the original ranking order for a few coin pairs that had been
# swapped (Cardano/Dogecoin, Polygon/Polkadot, Hedera/ImmutableX,
# Mantle/RenderToken).
#
# Price values in column D are plain numeric-looking text (e.g. "29.255.76",
# "0.7139") that must stay text, exactly as authored, rather than be
# reinterpreted by Excel as numbers (which would silently drop trailing
# zeros, switch to scientific notation, or fail outright for values using
# '.' as both grouping and decimal separators). Prefixing the value with a
# leading single-quote forces Excel to store it as text while keeping the
# displayed/read value free of the marker itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.255.76'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '''1.864.21'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''0.7139'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").Value = '''241.07'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").Value = '''0.07731'
$ws.Range("E8").Value = '  -1.03%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '''0.3084'
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("D10").Value = '''24.90'
$ws.Range("E10").Value = '  -0.62%  '
$ws.Range("D11").Value = '''0.08332'
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("D12").Value = '''1.874.21'
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''5.208'
$ws.Range("E13").Value = '  -1.70%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '''0.7145'
$ws.Range("E14").Value = '  -1.96%  '
$ws.Range("D15").Value = '''90.96'
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").Value = '''29.283.93'
$ws.Range("E16").Value = '  -0.51%  '
$ws.Range("D17").Value = '''5.964'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").Value = '''242.87'
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("D19").Value = '''0.000007823'
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("D20").Value = '''2.129.63'
$ws.Range("E20").Value = '  +0.50%  '
$ws.Range("D21").Value = '''13.17'
$ws.Range("E21").Value = '  -1.14%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").Value = '''7.907'
$ws.Range("E23").Value = '  -0.68%  '
$ws.Range("D24").Value = '''1.001'
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").Value = '''0.1603'
$ws.Range("E25").Value = '  +1.99%  '
$ws.Range("D26").Value = '''163.42'
$ws.Range("D27").Value = '''8.898'
$ws.Range("E27").Value = '  -1.69%  '
$ws.Range("D28").Value = '''18.57'
$ws.Range("E28").Value = '  +1.30%  '
$ws.Range("D29").Value = '''1.348'
$ws.Range("E29").Value = '  -1.43%  '
$ws.Range("D30").Value = '''1.499'
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("D31").Value = '''4.426'
$ws.Range("E31").Value = '  +0.90%  '
$ws.Range("D32").Value = '''4.260'
$ws.Range("E32").Value = '  +2.61%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''0.8414'
$ws.Range("E33").Value = '  +16.52%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '''0.05162'
$ws.Range("E34").Value = '  -2.30%  '
$ws.Range("D35").Value = '''1.931'
$ws.Range("E35").Value = '  -0.85%  '
$ws.Range("D36").Value = '''1.172'
$ws.Range("E36").Value = '  -2.60%  '
$ws.Range("E37").Value = '  +0.21%  '
$ws.Range("D38").Value = '''0.01855'
$ws.Range("E38").Value = '  -0.54%  '
$ws.Range("D39").Value = '''2.693'
$ws.Range("E39").Value = '  -1.10%  '
$ws.Range("D40").Value = '''1.165.76'
$ws.Range("E40").Value = '  -5.67%  '
$ws.Range("D41").Value = '''6.198'
$ws.Range("E41").Value = '  +1.53%  '
$ws.Range("D42").Value = '''0.8955'
$ws.Range("E42").Value = '  -1.42%  '
$ws.Range("D43").Value = '''72.86'
$ws.Range("E43").Value = '  -1.33%  '
$ws.Range("D44").Value = '''1.001'
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '''102.27'
$ws.Range("E45").Value = '  -1.41%  '
$ws.Range("D46").Value = '''2.026.39'
$ws.Range("E46").Value = '  +0.26%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''1.787'
$ws.Range("E47").Value = '  +1.25%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '''0.5169'
$ws.Range("E48").Value = '  -3.30%  '
$ws.Range("D49").Value = '''9.353'
$ws.Range("E49").Value = '  +0.63%  '
$ws.Range("D50").Value = '''1.002'
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("D51").Value = '''7.060'
$ws.Range("E51").Value = '  -0.38%  '
